# fix(OpportunityDataMapping): Fixed Issue with the VAT & Contact Mapping
# - Mapped `Contact Type` to `Contact` on Import
#
# This adds a second contact data row (row 4) to the "Contact" sheet,
# duplicating the formatting of the existing sample row (row 3) and
# changing the "No." value to 2 and the "Type" column to "Software"
# (instead of "Hardware"), exercising the corrected Contact-Type mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3's formatting (styles) and values into row 4 without
# creating new style entries, by round-tripping through Copy / PasteSpecial.
$ws.Range("A3:AA3").Copy()
$ws.Range("A4:AA4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A3:AA3").Copy()
$ws.Range("A4:AA4").PasteSpecial(-4163)   # xlPasteValues

# Row 3's auto-sized row height (driven by the wrapped "Street address"
# text) so row 4 matches visually.
$ws.Rows("4:4").RowHeight = 63.75

# New row-specific data: second contact ("No." = 2) with Contact Type
# mapped to "Software" rather than the "Hardware" value copied from row 3.
$ws.Range("A4").Value = 2
$ws.Range("M4").Value = "Software"
